$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.617.59'
$ws.Range("E2").Value = '  -4.52%  '
$ws.Range("D3").Value = '3.336.08'
$ws.Range("E3").Value = '  -1.44%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '574.09'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.36%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '180.48'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.66%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.629'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.18%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  -3.63%  '
$ws.Range("E10").Value = '  -1.57%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.403'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.07%  '
$ws.Range("D12").Value = '3.912.70'
$ws.Range("E12").Value = '  -1.56%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.136'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.61%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.01'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -6.03%  '
$ws.Range("D15").Value = '66.722.04'
$ws.Range("E15").Value = '  -4.33%  '
$ws.Range("E16").Value = '  -3.02%  '
$ws.Range("D17").Value = '3.343.02'
$ws.Range("E17").Value = '  -1.13%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '437.01'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.86%  '
$ws.Range("E19").Value = '  -2.38%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.53'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.59'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.88%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.53'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.43%  '
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("E24").Value = '  -0.99%  '
$ws.Range("E25").Value = '  -4.25%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.191'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.03'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.54%  '
$ws.Range("E28").Value = '  +0.10%  '
$ws.Range("E29").Value = '  -3.63%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.80'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.84%  '
$ws.Range("E31").Value = '  +0.01%  '
$ws.Range("E32").Value = '  -6.56%  '
$ws.Range("B33").Value = 'Aptos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.76'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.62%  '
$ws.Range("B34").Value = 'Fetch.AI'
$ws.Range("C34").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.23'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.65%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '163.20'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.78%  '
$ws.Range("E36").Value = '  -6.03%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '27.37'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.11%  '
$ws.Range("E38").Value = '  -8.24%  '
$ws.Range("D39").Value = '2.832.59'
$ws.Range("E39").Value = '  +3.04%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.796'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.28%  '
$ws.Range("E41").Value = '  -4.10%  '
$ws.Range("E42").Value = '  -5.59%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.18'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.34%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0667'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.25%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '24.42'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.40%  '
$ws.Range("E46").Value = '  -7.32%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '321.72'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.35%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0274'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.92%  '
$ws.Range("E49").Value = '  +0.56%  '
$ws.Range("E50").Value = '  -4.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.17'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.74%  '
